$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.164.18"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").Value = "2.254.80"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'307.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.91%  "

$ws.Range("D6").Value = "'98.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.80%  "

$ws.Range("D7").Value = "'0.575"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.79%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.56%  "

$ws.Range("D10").Value = "'35.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.50%  "

$ws.Range("D11").Value = "'0.0825"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "

$ws.Range("D12").Value = "'7.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.84%  "

$ws.Range("E13").Value = "  -1.80%  "

$ws.Range("D14").Value = "2.595.68"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.843"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.254.48"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "'13.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("D18").Value = "44.027.21"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("D19").Value = "'12.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.93%  "

$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").Value = "'6.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.32%  "

$ws.Range("D22").Value = "'65.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").Value = "'241.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.17%  "

$ws.Range("D24").Value = "'2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.97%  "

$ws.Range("E25").Value = "  -8.35%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'10.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").Value = "'37.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  -2.35%  "

$ws.Range("D30").Value = "'6.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.24%  "

$ws.Range("D31").Value = "'20.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").Value = "'157.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.00%  "

$ws.Range("D33").Value = "'3.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.09%  "

$ws.Range("D34").Value = "'0.0828"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.11%  "

$ws.Range("E35").Value = "  -1.61%  "

$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("E37").Value = "  -5.42%  "

$ws.Range("E38").Value = "  -4.50%  "

$ws.Range("D39").Value = "'15.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").Value = "'3.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.19%  "

$ws.Range("D41").Value = "'3.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.14%  "

$ws.Range("E42").Value = "  -3.94%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "1.757.95"
$ws.Range("E44").Value = "  -3.48%  "

$ws.Range("D45").Value = "'88.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.41%  "

$ws.Range("D46").Value = "'0.194"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.51%  "

$ws.Range("D47").Value = "'5.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("D48").Value = "'101.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'70.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.38%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.18%  "

$ws.Range("D51").Value = "'55.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.64%  "
